$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 518-519; existing rows 518:533 shift down to 520:535.
$ws.Rows("518:519").Insert()

# New row 518 - Naranja / Navel Late / Primera
$ws.Range("A518").Value = 10
$ws.Range("B518").Value = "Vega Modelo de Temuco"
$ws.Range("C518").Value = "La Araucanía"
$ws.Range("D518").Value = 44509
$ws.Range("E518").Value = 9
$ws.Range("F518").Value = "Fruta"
$ws.Range("G518").Value = 100102
$ws.Range("H518").Value = "Cítricos"
$ws.Range("I518").Value = 100102005
$ws.Range("J518").Value = "Naranja"
$ws.Range("K518").Value = "Navel Late"
$ws.Range("L518").Value = "Primera"
$ws.Range("M518").Value = 155
$ws.Range("N518").Value = 10000
$ws.Range("O518").Value = 10000
$ws.Range("P518").Value = 10000
$ws.Range("Q518").Value = "$/bandeja 15 kilos granel"
$ws.Range("R518").Value = "Región de O'Higgins"
$ws.Range("S518").Value = 667
$ws.Range("T518").Value = 15

# New row 519 - Naranja / Navel Late / Segunda
$ws.Range("A519").Value = 10
$ws.Range("B519").Value = "Vega Modelo de Temuco"
$ws.Range("C519").Value = "La Araucanía"
$ws.Range("D519").Value = 44509
$ws.Range("E519").Value = 9
$ws.Range("F519").Value = "Fruta"
$ws.Range("G519").Value = 100102
$ws.Range("H519").Value = "Cítricos"
$ws.Range("I519").Value = 100102005
$ws.Range("J519").Value = "Naranja"
$ws.Range("K519").Value = "Navel Late"
$ws.Range("L519").Value = "Segunda"
$ws.Range("M519").Value = 125
$ws.Range("N519").Value = 8000
$ws.Range("O519").Value = 8000
$ws.Range("P519").Value = 8000
$ws.Range("Q519").Value = "$/bandeja 15 kilos granel"
$ws.Range("R519").Value = "Región de O'Higgins"
$ws.Range("S519").Value = 533
$ws.Range("T519").Value = 15

$ws.Range("D518:D519").NumberFormat = "YYYY-MM-DD HH:MM:SS"
